# ERK fixed inverse colors bug
# Replace the semi-transparent accent-blue/orange rectangle fills on slide 3
# with their opaque "corrected" colors (the a:alpha children go away once the
# shape's solid fill is reassigned through Fill.ForeColor.RGB without
# touching Fill.Transparency).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$colorMap = @{
    "Rectangle 2"  = 8380415   # 4472C4 (50% alpha) -> FFDF7F
    "Rectangle 14" = 8571647   # 235888 (50% alpha) -> FFCA82
    "Rectangle 16" = 10930159  # 2C70AE (75% alpha) -> EFC7A6
    "Rectangle 18" = 10533852  # BA7741 (50% alpha) -> DCBBA0
    "Rectangle 20" = 12751969  # E08F4E (50% alpha) -> 6194C2
    "Rectangle 22" = 12823441  # FFB858 (75% alpha) -> 91ABC3
    "Rectangle 24" = 14792865  # FFC000 (50% alpha) -> A1B8E1
}

foreach ($name in $colorMap.Keys) {
    $shape = $s.Shapes.Item($name)
    $shape.Fill.ForeColor.RGB = $colorMap[$name]
}
